$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last data row (44) currently holds the "NA" page-number value in
# column C. The script re-run produced a new day's result (2025-04-10)
# that is also "Rien ne nous concerne aujourd'hui !" with "NA" for the
# page number, so that new row is appended as row 45, and row 44's C
# value is cleared (matching the other "nothing to report" rows whose
# C column is blank).

# 1) Append the new row 45 with the latest scraped result.
#    Column A holds a date-like string (e.g. "2025-04-09" on row 44) but
#    must stay a literal text value, not get auto-converted to a date
#    serial number. Force text via NumberFormat, then restore the
#    original (unstyled) formatting by pasting A44's formats over it so
#    no stray style is introduced.
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = "2025-04-10"
$ws.Cells.Item(44, 1).Copy()
$ws.Cells.Item(45, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(45, 2).Value = "Rien ne nous concerne aujourd'hui !"
$ws.Cells.Item(45, 3).Value = "NA"
$ws.Cells.Item(45, 4).Value = 1

# 2) Clear the old "NA" from row 44's C column now that it has moved to
#    row 45 - row 44 becomes a normal blank-C "nothing to report" row.
$ws.Cells.Item(44, 3).ClearContents()
